$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row total (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row correct-marks (B12): 63 -> 105
$ws.Range("B12").Value = 105

# Update corresponding "corr/total" label (E12): "63/84" -> "105/140"
$ws.Range("E12").Value = "105/140"
